$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New identifiers that replace the old ones throughout the workbook.
# ---------------------------------------------------------------------------
$oldUuid1 = "9606c80f-3135-4c57-8498-146ae6416c79"
$newUuid1 = "4f2cfc93-10be-49d1-be8c-2e3bd74551aa"
$oldUuid2 = "a01eb00f-f7c8-4c83-b44f-fa64fe24f78b"
$newUuid2 = "ffff6e5c7937-abe1-4561-b50b-14bbb41cdd4a"

$newStatus = "Ready for handoff"
$newHoDate = "2016-09-07 09:37:40"
$newHandoffDate = "2016-09-07 09:37:34"
$newHandbackDate = "0001-01-01 00:00:00"

$newHash = "9114210121402a8c80bb89005014ed5357fa7d5f"
$newZhXlf = "$newUuid1.$newHash.zh-cn.xlf"
$newDeXlf = "$newUuid1.$newHash.de-de.xlf"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newUuid1.md"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Range("A3").Value = "$newUuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newUuid2.md"
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $newHoDate

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newUuid1.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$newUuid2.md"
    }
}

$wsOverview.Columns.Item(5).ColumnWidth = 16.38265482584637
$wsOverview.Columns.Item(6).ColumnWidth = 16.38265482584637

# ---------------------------------------------------------------------------
# zh-cn / de-de sheets share the same layout - apply identical edits to both,
# parameterised by the language-specific xlf name + HO date.
# ---------------------------------------------------------------------------
function Update-LangSheet($ws, $xlfName, $hoDate) {
    $ws.Range("A2").Value = "$newUuid1.md"
    $ws.Range("C2").Value = $newStatus
    $ws.Range("G2").Value = $xlfName
    $ws.Range("H2").Value = $newHandoffDate
    $ws.Range("I2").Value = ""
    $ws.Range("J2").Value = ""
    $ws.Range("K2").Value = $newHandbackDate

    $ws.Range("A3").Value = "$newUuid2.md"
    $ws.Range("C3").Value = $newStatus
    $ws.Range("F3").Value = "True"
    $ws.Range("G3").Value = $xlfName
    $ws.Range("H3").Value = $newHandoffDate
    $ws.Range("I3").Value = ""
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = $newHandbackDate

    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$I$2' -or $addr -eq '$I$3') {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }

    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$A$2') {
            $hl.TextToDisplay = "$newUuid1.md"
        } elseif ($addr -eq '$A$3') {
            $hl.TextToDisplay = "$newUuid2.md"
        }
    }

    $ws.Columns.Item(3).ColumnWidth = 16.38265482584637
    $ws.Columns.Item(9).ColumnWidth = 17.817272004627068
    $ws.Columns.Item(10).ColumnWidth = 20.872143700009268
}

Update-LangSheet $wsZh $newZhXlf $newHoDate
Update-LangSheet $wsDe $newDeXlf $newHoDate
